$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: toggle values in column D for several rows ---
# (rows < 26, so unaffected by the later row deletions)
$ws.Range("D6").Value = -14.2
$ws.Range("D8").ClearContents()
$ws.Range("D12").Value = -14.1
$ws.Range("D14").ClearContents()
$ws.Range("D17").Value = -14.7
$ws.Range("D18").Value = -15.2
$ws.Range("D19").ClearContents()
$ws.Range("D20").ClearContents()
$ws.Range("D23").Value = -13.9

# --- Step 2: remove the "RM 232" row (row 26) entirely ---
$ws.Rows.Item(26).Delete()

# --- Step 3: remove the "SC 92" row (now shifted up to row 27) entirely ---
$ws.Rows.Item(27).Delete()

# --- Step 4: apply remaining value changes on the now-shifted rows ---
# Row 27 = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("D27").ClearContents()

# Row 28 = "SC 105"
$ws.Range("B28").ClearContents()

# Row 29 = "SC 119"
$ws.Range("B29").ClearContents()

# Row 30 = "SC 120"
$ws.Range("B30").Value = -19.7

# Row 32 = "SC 193"
$ws.Range("B32").ClearContents()
